# Updates the cryptos list (Price / Volume(1h) columns) to new scraped values.
# Note: some Price values (single-decimal-point numbers, e.g. "216.35") are
# prefixed with a leading apostrophe so Excel stores them as text (matching
# the original cells, which are all text) instead of auto-converting them to
# numbers. Values that already contain two dots (thousands separators, e.g.
# "26.042.48") are never auto-parsed as numbers so they need no prefix.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.042.48"
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("D3").Value = "1.644.72"
$ws.Range("E3").Value = "  +0.88%  "
$ws.Range("E4").Value = "  +0.71%  "
$ws.Range("D5").Value = "'216.35"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("E6").Value = "  +0.97%  "
$ws.Range("E7").Value = "  +0.63%  "
$ws.Range("E8").Value = "  +0.54%  "
$ws.Range("E9").Value = "  +1.18%  "
$ws.Range("D10").Value = "'19.65"
$ws.Range("E10").Value = "  +0.03%  "
$ws.Range("E11").Value = "  +1.23%  "
$ws.Range("D12").Value = "1.873.38"
$ws.Range("E12").Value = "  +0.93%  "
$ws.Range("D14").Value = "1.650.51"
$ws.Range("E14").Value = "  +1.33%  "
$ws.Range("E15").Value = "  -0.08%  "
$ws.Range("E16").Value = "  +0.94%  "
$ws.Range("D17").Value = "'63.23"
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("D18").Value = "26.046.56"
$ws.Range("E18").Value = "  +0.70%  "
$ws.Range("E19").Value = "  +0.66%  "
$ws.Range("D20").Value = "'193.09"
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("E21").Value = "  -0.69%  "
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("E24").Value = "  +5.18%  "
$ws.Range("E25").Value = "  +0.53%  "
$ws.Range("D26").Value = "'144.45"
$ws.Range("E26").Value = "  +1.27%  "
$ws.Range("E27").Value = "  +0.87%  "
$ws.Range("D28").Value = "'6.91"
$ws.Range("E28").Value = "  +0.71%  "
$ws.Range("E29").Value = "  +0.63%  "
$ws.Range("E30").Value = "  +1.10%  "
$ws.Range("E31").Value = "  +0.31%  "
$ws.Range("E32").Value = "  -0.38%  "
$ws.Range("D33").Value = "'3.26"
$ws.Range("E33").Value = "  +1.23%  "
$ws.Range("E34").Value = "  -2.88%  "
$ws.Range("E35").Value = "  +2.41%  "
$ws.Range("E36").Value = "  +0.58%  "
$ws.Range("D37").Value = "1.131.95"
$ws.Range("E37").Value = "  -0.37%  "
$ws.Range("D38").Value = "'0.540"
$ws.Range("E39").Value = "  +0.19%  "
$ws.Range("E40").Value = "  +0.81%  "
$ws.Range("D41").Value = "'5.52"
$ws.Range("E41").Value = "  +0.83%  "
$ws.Range("D42").Value = "'99.45"
$ws.Range("E42").Value = "  +0.30%  "
$ws.Range("E43").Value = "  -0.52%  "
$ws.Range("D44").Value = "1.782.67"
$ws.Range("E44").Value = "  +0.96%  "
$ws.Range("E45").Value = "  +3.83%  "
$ws.Range("D46").Value = "'56.59"
$ws.Range("E46").Value = "  +0.92%  "
$ws.Range("E47").Value = "  -0.29%  "
$ws.Range("E48").Value = "  +0.28%  "
$ws.Range("D49").Value = "'7.71"
$ws.Range("E49").Value = "  +1.08%  "
$ws.Range("E50").Value = "  +0.43%  "
$ws.Range("E51").Value = "  -0.30%  "
